$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header casing: "Symbol" -> "symbol"
$ws.Range("A1").Value = "symbol"

# The B:D columns only ever held empty, pre-formatted placeholder cells
# (rows 2-9) left over from when the sheet was first built; clear them out
# entirely so they no longer exist in the saved file.
$ws.Range("B2:D9").Clear()

# Column A's width was sized for the old (wider) sample data; re-fit it to
# the current contents now that the sheet only needs to show the symbol list.
$ws.Columns("A").AutoFit()

# Match the saved selection/active cell.
[void]$ws.Range("A2").Select()
